# Isis.xlsx BO workbook update:
#  - add a new YouTube hyperlink row on Blad1 (sheet1) at A19
#  - add a new YouTube hyperlink row on Feuil1 (sheet2) at A27 (with tooltip)
#  - add a new plain-text (Discord/video) link row on Feuil2 (sheet3) at A13
#
# Shared-string append order matters (it determines the resulting sst index),
# so cells are populated in the same order the target workbook expects:
#   1) Feuil1!A27  -> https://www.youtube.com/watch?v=h-2EuniduXM&t=&ab_channel=DeitiesofDeath
#   2) Blad1!A19   -> https://youtu.be/h-2EuniduXM?si=DzI9yv7mW9gFL4KP&t=510
#   3) Feuil2!A13  -> https://youtu.be/h-2EuniduXM?si=BsCb5yP7MBGUuLad&t=832

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Blad1")
$ws2 = $wb.Worksheets.Item("Feuil1")
$ws3 = $wb.Worksheets.Item("Feuil2")

# 1) Feuil1 (sheet2) row 27 - hyperlink with an explicit tooltip/screentip
$url2 = "https://www.youtube.com/watch?v=h-2EuniduXM&t=&ab_channel=DeitiesofDeath"
$ws2.Range("A27").Value = $url2
$ws2.Hyperlinks.Add($ws2.Range("A27"), $url2, "", $url2)

# 2) Blad1 (sheet1) row 19 - hyperlink, no tooltip
$url1 = "https://youtu.be/h-2EuniduXM?si=DzI9yv7mW9gFL4KP&t=510"
$ws1.Range("A19").Value = $url1
$ws1.Hyperlinks.Add($ws1.Range("A19"), $url1)

# 3) Feuil2 (sheet3) row 13 - plain text, no hyperlink
$url3 = "https://youtu.be/h-2EuniduXM?si=BsCb5yP7MBGUuLad&t=832"
$ws3.Range("A13").Value = $url3

# Restore the selections to the newly-added cells on each sheet
$ws1.Range("A19").Select()
$ws2.Range("A27").Select()
$ws3.Range("A13").Select()
